# Insert a new data row before the current row 259, shifting rows 259-310
# down to 260-311, and populate the newly inserted row 259 with the new
# record's values (dated 2022-03-17 / serial 44637).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 259; this pushes the existing rows 259..310
# down to 260..311 and extends the sheet dimension to A1:R311.
$ws.Rows.Item(259).Insert()

# Populate the new row 259 with the new record.
$ws.Range("A259").Value = 5
$ws.Range("B259").Value = "Macroferia Regional de Talca"
$ws.Range("C259").Value = "Maule"
$ws.Range("D259").Value = 44637
$ws.Range("E259").Value = 7
$ws.Range("F259").Value = 100112032
$ws.Range("G259").Value = "Zapallo italiano"
$ws.Range("H259").Value = "Sin especificar"
$ws.Range("I259").Value = "Primera"
$ws.Range("J259").Value = 500
$ws.Range("K259").Value = 8000
$ws.Range("L259").Value = 8000
$ws.Range("M259").Value = 8000
$ws.Range("N259").Value = "$/caja 50 unidades"
$ws.Range("O259").Value = "Región del Maule"
$ws.Range("P259").Value = 160
$ws.Range("Q259").Value = 50
$ws.Range("R259").Value = "Hortaliza"
